$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (column C) and montant_total (column D) values
# for rows with newly published 2020-09-27 data. Values are kept as text
# (NumberFormat "@") to match the source workbook's inline-string cell type
# and to preserve exact decimal formatting (e.g. trailing zeros).

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1379"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "3416243.70"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "1110"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4445207.81"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "750"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2677026.29"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "39"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "170144.45"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "223"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "685665.52"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "399"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1560611.18"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "186"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "771545.34"

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "248"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "665219.00"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "527"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2021649.55"

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "154"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "551601.67"

$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "200"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "536288.00"

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "358"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1488305.26"

$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "171"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "633289.45"

$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "11"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "26500.00"

$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "336"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "889341.89"

$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "656"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3194383.47"

$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "435"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1610094.99"

$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "17"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "66860.00"

$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "19"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "61932.00"

$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "484"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1260613.65"

$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "244"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "817750.04"

$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "243"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "795974.14"

$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "8"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19000.00"

$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "21"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "53623.85"

$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "265"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "756399.15"

$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "125"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "722333.98"

$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "178"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "710032.25"

$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "7"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23519.00"

$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "522"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1466838.34"

$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "751"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3502717.57"

$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "525"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2144483.62"

$ws.Range("C53").NumberFormat = "@"
$ws.Range("C53").Value = "6400"
$ws.Range("D53").NumberFormat = "@"
$ws.Range("D53").Value = "15683066.28"

$ws.Range("C56").NumberFormat = "@"
$ws.Range("C56").Value = "31"
$ws.Range("D56").NumberFormat = "@"
$ws.Range("D56").Value = "205600.00"

$ws.Range("C57").NumberFormat = "@"
$ws.Range("C57").Value = "4874"
$ws.Range("D57").NumberFormat = "@"
$ws.Range("D57").Value = "19287073.30"

$ws.Range("C59").NumberFormat = "@"
$ws.Range("C59").Value = "4964"
$ws.Range("D59").NumberFormat = "@"
$ws.Range("D59").Value = "16880511.05"

$ws.Range("C61").NumberFormat = "@"
$ws.Range("C61").Value = "106"
$ws.Range("D61").NumberFormat = "@"
$ws.Range("D61").Value = "409614.47"

$ws.Range("C63").NumberFormat = "@"
$ws.Range("C63").Value = "176"
$ws.Range("D63").NumberFormat = "@"
$ws.Range("D63").Value = "551954.21"

$ws.Range("C70").NumberFormat = "@"
$ws.Range("C70").Value = "50"
$ws.Range("D70").NumberFormat = "@"
$ws.Range("D70").Value = "171460.00"

$ws.Range("C71").NumberFormat = "@"
$ws.Range("C71").Value = "25"
$ws.Range("D71").NumberFormat = "@"
$ws.Range("D71").Value = "98866.05"

$ws.Range("C78").NumberFormat = "@"
$ws.Range("C78").Value = "419"
$ws.Range("D78").NumberFormat = "@"
$ws.Range("D78").Value = "1110329.60"

$ws.Range("C80").NumberFormat = "@"
$ws.Range("C80").Value = "1007"
$ws.Range("D80").NumberFormat = "@"
$ws.Range("D80").Value = "4194170.88"

$ws.Range("C81").NumberFormat = "@"
$ws.Range("C81").Value = "557"
$ws.Range("D81").NumberFormat = "@"
$ws.Range("D81").Value = "2209716.12"

$ws.Range("C82").NumberFormat = "@"
$ws.Range("C82").Value = "42"
$ws.Range("D82").NumberFormat = "@"
$ws.Range("D82").Value = "132078.00"

$ws.Range("C83").NumberFormat = "@"
$ws.Range("C83").Value = "38"
$ws.Range("D83").NumberFormat = "@"
$ws.Range("D83").Value = "168180.27"

$ws.Range("C84").NumberFormat = "@"
$ws.Range("C84").Value = "671"
$ws.Range("D84").NumberFormat = "@"
$ws.Range("D84").Value = "1692707.33"

$ws.Range("C87").NumberFormat = "@"
$ws.Range("C87").Value = "1015"
$ws.Range("D87").NumberFormat = "@"
$ws.Range("D87").Value = "3824759.84"

$ws.Range("C88").NumberFormat = "@"
$ws.Range("C88").Value = "727"
$ws.Range("D88").NumberFormat = "@"
$ws.Range("D88").Value = "2325622.16"

$ws.Range("C91").NumberFormat = "@"
$ws.Range("C91").Value = "237"
$ws.Range("D91").NumberFormat = "@"
$ws.Range("D91").Value = "563350.00"

$ws.Range("C93").NumberFormat = "@"
$ws.Range("C93").Value = "540"
$ws.Range("D93").NumberFormat = "@"
$ws.Range("D93").Value = "2060616.16"

$ws.Range("C94").NumberFormat = "@"
$ws.Range("C94").Value = "202"
$ws.Range("D94").NumberFormat = "@"
$ws.Range("D94").Value = "634423.09"

$ws.Range("C97").NumberFormat = "@"
$ws.Range("C97").Value = "929"
$ws.Range("D97").NumberFormat = "@"
$ws.Range("D97").Value = "2305321.17"

$ws.Range("C100").NumberFormat = "@"
$ws.Range("C100").Value = "1252"
$ws.Range("D100").NumberFormat = "@"
$ws.Range("D100").Value = "4574359.41"

$ws.Range("C102").NumberFormat = "@"
$ws.Range("C102").Value = "1198"
$ws.Range("D102").NumberFormat = "@"
$ws.Range("D102").Value = "4144473.98"
